# Auto-generated: rewrite rows 2-29 (BSL data) to match target snapshot (v15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 28,58
for ($r = 0; $r -lt $data.GetLength(0); $r++) {
  for ($c = 0; $c -lt $data.GetLength(1); $c++) {
    $data[$r,$c] = 0
  }
}

# Row 2 (array row 0)
$data[0,0] = 2   # A2
$data[0,1] = 'Biltin'   # B2
$data[0,2] = 'Biltin 20mg Tablet 30''s'   # C2
$data[0,3] = '30''s'   # D2
$data[0,15] = 312   # P2
$data[0,17] = 312   # R2
$data[0,54] = 337.33   # BC2

# Row 3 (array row 1)
$data[1,0] = 4   # A3
$data[1,1] = 'Desodin'   # B3
$data[1,2] = 'Desodin 60ml Syrup'   # C3
$data[1,3] = '60 ml'   # D3
$data[1,54] = 18.74   # BC3

# Row 4 (array row 2)
$data[2,0] = 5   # A4
$data[2,1] = 'Dinafex'   # B4
$data[2,2] = 'Dinafex 180mg Tablet'   # C4
$data[2,3] = '30''s'   # D4
$data[2,54] = 224.89   # BC4

# Row 5 (array row 3)
$data[3,0] = 5   # A5
$data[3,1] = 'Dinafex'   # B5
$data[3,2] = 'Dinafex 120mg Tablet'   # C5
$data[3,3] = '30''s'   # D5
$data[3,54] = 179.91   # BC5

# Row 6 (array row 4)
$data[4,0] = 5   # A6
$data[4,1] = 'Dinafex'   # B6
$data[4,2] = 'Dinafex 60mg Tablet'   # C6
$data[4,3] = '30''s'   # D6
$data[4,54] = 78.70999999999999   # BC6

# Row 7 (array row 5)
$data[5,0] = 6   # A7
$data[5,1] = 'Dorenta'   # B7
$data[5,2] = 'Dorenta 50mg Tablet'   # C7
$data[5,3] = '50''s'   # D7
$data[5,54] = 93.70999999999999   # BC7

# Row 8 (array row 6)
$data[6,0] = 7   # A8
$data[6,1] = 'Etorix'   # B8
$data[6,2] = 'Etorix 90mg Tablet'   # C8
$data[6,3] = '30''s'   # D8
$data[6,54] = 269.87   # BC8

# Row 9 (array row 7)
$data[7,0] = 7   # A9
$data[7,1] = 'Etorix'   # B9
$data[7,2] = 'Etorix 120mg Tablet'   # C9
$data[7,3] = '20''s'   # D9
$data[7,54] = 209.9   # BC9

# Row 10 (array row 8)
$data[8,0] = 7   # A10
$data[8,1] = 'Etorix'   # B10
$data[8,2] = 'Etorix 60mg Tablet - 40''s'   # C10
$data[8,3] = '40''s'   # D10
$data[8,54] = 209.9   # BC10

# Row 11 (array row 9)
$data[9,0] = 8   # A11
$data[9,1] = 'Fenobac'   # B11
$data[9,2] = 'Fenobac 100ml Syrup'   # C11
$data[9,3] = '100ml'   # D11
$data[9,54] = 74.95999999999999   # BC11

# Row 12 (array row 10)
$data[10,0] = 9   # A12
$data[10,1] = 'Flucloxin'   # B12
$data[10,2] = 'Flucloxin 500mg Capsule'   # C12
$data[10,3] = '30 ''s'   # D12
$data[10,54] = 237.74   # BC12

# Row 13 (array row 11)
$data[11,0] = 9   # A13
$data[11,1] = 'Flucloxin'   # B13
$data[11,2] = 'Flucloxin 500mg Capsule - 36''s'   # C13
$data[11,3] = '36 ''s'   # D13
$data[11,54] = 284.21   # BC13

# Row 14 (array row 12)
$data[12,0] = 10   # A14
$data[12,1] = 'Geminox'   # B14
$data[12,2] = 'Geminox 320mg Tablet - 8''s'   # C14
$data[12,3] = '8 ''s'   # D14
$data[12,6] = 310   # G14
$data[12,7] = 258   # H14
$data[12,54] = 389.8   # BC14

# Row 15 (array row 13)
$data[13,0] = 11   # A15
$data[13,1] = 'Ketonic'   # B15
$data[13,2] = 'Ketonic 30mg Injection'   # C15
$data[13,3] = '5 ''s'   # D15
$data[13,54] = 206.77   # BC15

# Row 16 (array row 14)
$data[14,0] = 11   # A16
$data[14,1] = 'Ketonic'   # B16
$data[14,2] = 'Ketonic 10mg Tablet'   # C16
$data[14,3] = '20''s'   # D16
$data[14,54] = 150.38   # BC16

# Row 17 (array row 15)
$data[15,0] = 11   # A17
$data[15,1] = 'Ketonic'   # B17
$data[15,2] = 'Ketonic 30mg IM/IV Injection - 4''s'   # C17
$data[15,3] = '4''s'   # D17
$data[15,54] = 165.41   # BC17

# Row 18 (array row 16)
$data[16,0] = 12   # A18
$data[16,1] = 'Kynol'   # B18
$data[16,2] = 'Kynol TR 200mg Capsule'   # C18
$data[16,3] = '30 ''s'   # D18
$data[16,54] = 224.89   # BC18

# Row 19 (array row 17)
$data[17,0] = 12   # A19
$data[17,1] = 'Kynol'   # B19
$data[17,2] = 'Kynol D 25mg Tablet'   # C19
$data[17,3] = '60 ''s'   # D19
$data[17,54] = 180.45   # BC19

# Row 20 (array row 18)
$data[18,0] = 12   # A20
$data[18,1] = 'Kynol'   # B20
$data[18,2] = 'Kynol TR 100mg Capsule'   # C20
$data[18,3] = '50 ''s'   # D20
$data[18,54] = 262.37   # BC20

# Row 21 (array row 19)
$data[19,0] = 17   # A21
$data[19,1] = 'Naprox'   # B21
$data[19,2] = 'Naprox Plus 500mg Tablet - 30''s'   # C21
$data[19,3] = '30 ''s'   # D21
$data[19,54] = 224.89   # BC21

# Row 22 (array row 20)
$data[20,0] = 19   # A22
$data[20,1] = 'Oradin'   # B22
$data[20,2] = 'Oradin Plus Tablet - 40''s'   # C22
$data[20,3] = '40 ''s'   # D22
$data[20,54] = 209.9   # BC22

# Row 23 (array row 21)
$data[21,0] = 20   # A23
$data[21,1] = 'Osticare'   # B23
$data[21,2] = 'Osticare Tablet 24''s'   # C23
$data[21,3] = '24''s'   # D23
$data[21,54] = 215.89   # BC23

# Row 24 (array row 22)
$data[22,0] = 23   # A24
$data[22,1] = 'Rupaday'   # B24
$data[22,2] = 'Rupaday Oral Solution 60ml'   # C24
$data[22,3] = '1''s'   # D24
$data[22,6] = 1458   # G24
$data[22,7] = 1215   # H24
$data[22,11] = 5   # L24
$data[22,13] = 1458   # N24
$data[22,54] = 56.22   # BC24

# Row 25 (array row 23)
$data[23,0] = 24   # A25
$data[23,1] = 'Sk-Mox'   # B25
$data[23,2] = 'Sk-Mox 500mg Capsule'   # C25
$data[23,3] = '48 ''s'   # D25
$data[23,8] = 14   # I25
$data[23,11] = 21   # L25
$data[23,14] = 173   # O25
$data[23,20] = 173   # U25
$data[23,38] = 21   # AM25
$data[23,47] = 152   # AV25
$data[23,52] = 150   # BA25
$data[23,53] = 32869   # BB25
$data[23,54] = 219.13   # BC25

# Row 26 (array row 24)
$data[24,0] = 35   # A26
$data[24,1] = 'Zithrox'   # B26
$data[24,2] = 'Zithrox 15ml Suspension'   # C26
$data[24,3] = '15 ml'   # D26
$data[24,54] = 71.95999999999999   # BC26

# Row 27 (array row 25)
$data[25,0] = 35   # A27
$data[25,1] = 'Zithrox'   # B27
$data[25,2] = 'Zithrox 500mg Tablet'   # C27
$data[25,3] = '6 ''s'   # D27
$data[25,54] = 136.83   # BC27

# Row 28 (array row 26)
$data[26,0] = 35   # A28
$data[26,1] = 'Zithrox'   # B28
$data[26,2] = 'Zithrox 30ml Dry Suspension'   # C28
$data[26,3] = '30ml'   # D28
$data[26,54] = 97.45   # BC28

# Row 29 (array row 27)
$data[27,0] = 35   # A29
$data[27,1] = 'Zithrox'   # B29
$data[27,2] = 'Zithrox 250mg Tablet - 6''s'   # C29
$data[27,3] = '6''s'   # D29
$data[27,54] = 89.95999999999999   # BC29

$ws.Range("A2:BF29").Value = $data